# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps for the 1255e723-... handback row on both the zh-cn and de-de
# sheets, as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-18 07:43:09"
$wsZhCn.Range("G4").Value = "2016-02-18 07:43:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-18 07:43:21"
$wsDeDe.Range("G4").Value = "2016-02-18 07:44:20"
